$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Assumptions sheet: add a "Current Share Price" input row, and remove
# the old equity-value/enterprise-value build-out (rows 9-14) that has
# been relocated to the Results sheet.
# ---------------------------------------------------------------------
$assump = $wb.Worksheets.Item("Assumptions")

# New input: Current Share Price = 5.75, formatted like the other
# currency inputs (copy format from B9, which still carries style s="4").
$assump.Range("A8").Value = "Current Share Price"
$assump.Range("B8").Value = 5.75
$assump.Range("B9").Copy()
$assump.Range("B8").PasteSpecial(-4122)

# Clear out the old derived rows 9-14 (formulas + labels), leaving the
# now-empty, but still formatted, cells behind -- row 13 had no special
# formatting, so clearing both its cells drops it from the sheet.
$assump.Range("A9:A14").ClearContents()
$assump.Range("B9:B14").ClearContents()

$assump.Range("F7").Select()

# ---------------------------------------------------------------------
# Results sheet: append the upside-vs-current-price calculations.
# ---------------------------------------------------------------------
$results = $wb.Worksheets.Item("Results")

$results.Range("A9").Value = "absolute_upside_per_share"
$results.Range("A8").Copy()
$results.Range("A9").PasteSpecial(-4122)

$results.Range("B9").Formula = "=B8-Assumptions!B8"
$assump.Range("B12").Copy()
$results.Range("B9").PasteSpecial(-4122)

$results.Range("A10").Value = "percent_upside_per_share"
$results.Range("A8").Copy()
$results.Range("A10").PasteSpecial(-4122)

$results.Range("B10").Formula = "=B9/Assumptions!B8"
$results.Range("B10").Style = "Percent"

$results.Range("F14").Select()
